# The invoice playground example is a tiny templated workbook: columns D/E
# on rows 2-3 still hold unresolved "[[items.date]]" / "[[items.missingProp]]"
# placeholders. This commit re-renders the example output with the
# placeholders substituted by their resolved values - row 2's item date
# resolves to 2025-01-01, row 3's item date resolves to 2025-01-02 (a value
# distinct from row 2, so it becomes a new shared string), and the missing
# property resolves to 0 for both rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Item date: 2025-01-01"
$ws.Range("E2").Value = "Missing: 0"

$ws.Range("D3").Value = "Item date: 2025-01-02"
$ws.Range("E3").Value = "Missing: 0"
